$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.090.24"
$ws.Range("E2").Value = "  -2.95%  "

$ws.Range("D3").Value = "1.652.80"
$ws.Range("E3").Value = "  -4.98%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.93%  "

$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4803"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.86%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2624"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.05982"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07105"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("D11").Value = "1.654.09"
$ws.Range("E11").Value = "  -4.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6209"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.602"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "73.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.51%  "

$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("D18").Value = "25.096.29"
$ws.Range("E18").Value = "  -3.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006541"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.427"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.62%  "

$ws.Range("D22").Value = "1.870.92"
$ws.Range("E22").Value = "  -4.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.476"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.280"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "133.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.390"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.45%  "

$ws.Range("E28").Value = "  -3.94%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "101.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.818"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07895"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.523"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04598"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.606"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9445"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5844"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.619"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01539"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8433"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.06%  "

$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.835"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3705"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.836"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1122"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.063"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05149"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3335"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.33%  "

